$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price updates (column D) scraped on Thu Dec 15 05:34:52 UTC 2022

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "264.85"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "22.71"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "6.283"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.06142"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "6.676"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.344"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.8293"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08226"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03424"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.03098"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.09248"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.913"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.001704"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.04884"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.006214"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.005269"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.291"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04613"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006945"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1137"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.003401"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01070"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00006164"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.7782"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.1950"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00002101"

$wb.Save()
